# Apply updated dSF (column F) values as part of "repull data, push all
# data, mean calculation".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 1
    5  = 0
    6  = 1
    11 = 1
    12 = -1
    16 = -3
    17 = 3
    18 = -3
    19 = -1
    26 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
